# Add a "Non-Atomix" baseline column to the "Atomix" sheet's perf table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Atomix")

# Insert a new column before the current column B (BASELINE), shifting
# BASELINE -> C and SC-SC-OPT -> D.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Non-Atomix"
$ws.Range("B1").Font.Bold = $true

# New "Non-Atomix" baseline values (rows 2..18).
$ws.Range("B2").Value = 141.738
$ws.Range("B3").Value = 83.485799999999998
$ws.Range("B4").Value = 130.785
$ws.Range("B5").Value = 234.142
$ws.Range("B6").Value = 50.476500000000001
$ws.Range("B7").Value = 60.640500000000003
$ws.Range("B8").Value = 99.204300000000003
$ws.Range("B9").Value = 221.26400000000001
$ws.Range("B10").Value = 53.211399999999998
$ws.Range("B11").Value = 119.5
$ws.Range("B12").Value = 50.630499999999998
$ws.Range("B13").Value = 197.06899999999999
$ws.Range("B14").Value = 112.996
$ws.Range("B15").Value = 230.12200000000001
$ws.Range("B16").Value = 67.674700000000001
$ws.Range("B17").Value = 84.192999999999998
$ws.Range("B18").Value = 234.87899999999999

# The row-label column is now bold as well (matches the rest of the header row).
$ws.Range("A2:A18").Font.Bold = $true

# Column widths for the (now 4) populated columns.
$ws.Range("B1").ColumnWidth = 13
$ws.Range("C1").ColumnWidth = 11.666666666666668
$ws.Range("D1").ColumnWidth = 12.833333333333332

# Match the author's final selection on the sheet.
[void]$ws.Range("B26").Select()
